$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on numeric-looking Price (D) values so Excel
# does not auto-convert them to numbers, matching the original inlineStr text cells.
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "63.212.26"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "3.051.45"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "588.99"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "152.68"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "3.055.47"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").Value = "5.89"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("D14").Value = "36.40"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "3.563.26"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "7.17"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "63.253.67"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "3.057.71"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "481.95"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "81.90"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "12.77"
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("D29").Value = "7.43"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.67"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "2.22"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "27.18"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "1.06"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0817"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("D37").Value = "6.01"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  -6.37%  "
$ws.Range("D39").Value = "2.20"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "9.29"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "50.50"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "436.99"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("D43").Value = "0.286"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "0.0361"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "2.814.57"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "39.28"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  -2.58%  "
